$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("shard")

# Swap the contents of rows 28 and 29 (columns A-C; D stays 21 in both).
$row28 = @($ws.Range("A28").Value2, $ws.Range("B28").Value2, $ws.Range("C28").Value2)
$row29 = @($ws.Range("A29").Value2, $ws.Range("B29").Value2, $ws.Range("C29").Value2)

$ws.Range("A28").Value = $row29[0]
$ws.Range("B28").Value = $row29[1]
$ws.Range("C28").Value = $row29[2]

$ws.Range("A29").Value = $row28[0]
$ws.Range("B29").Value = $row28[1]
$ws.Range("C29").Value = $row28[2]

# Both rows lose their explicit 12pt custom row height (back to default/auto).
$ws.Rows(28).AutoFit()
$ws.Rows(29).AutoFit()

# Update the active selection to match the final user action.
$ws.Range("B43").Select()
